# Generate Report for Handoff
# Updates the handoff timestamps and sets the "Priority" column (E) to "ht"
# for the Ready-for-handoff rows (7,8,9,10,11,14) on the zh-cn / de-de
# sheets, and bumps the corresponding "Latest HO/Handoff" timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 11, 14)

foreach ($r in $rows) {
    # Overview!G<r> - "Latest HO Xliff Generate Date" (text, not a real date)
    $overview.Range("G$r").Value = "2016-09-03 04:24:48"

    # de-de!H<r> - "Latest Handoff Datetime" (originally the same text as
    # Overview!G<r>, so it moves to the same new value)
    $dede.Range("H$r").Value = "2016-09-03 04:24:48"

    # zh-cn!H<r> - "Latest Handoff Datetime"
    $zhcn.Range("H$r").Value = "2016-09-03 04:24:43"

    # zh-cn!E<r> and de-de!E<r> - "Priority" column, empty -> "ht"
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
